# faturamento_diario.xlsx — "atualizei dados bibi e add"
#
# A new daily-revenue row (July / "07/2025", day 7) is inserted right after
# the existing July block (rows 2-7), pushing the June/May/April blocks
# (previously rows 8-98) down by one row to rows 9-99.
#
# Insert a fresh row at row 8 (shifts everything below it down one row,
# exactly like the diff), then fill it with the new July-day-7 data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(8).Insert()

$ws.Range("A8").Value = 7
$ws.Range("B8").Value = 19917.67
$ws.Range("C8").Value = 7
$ws.Range("D8").Value = 2025
$ws.Range("E8").Value = "07/2025"
